# feat: add 2022-Q1 data
#
# 1) Insert a new sheet "2022-Q1" right before "总计" with the fund-holding
#    detail table (same layout/header/styling as the other quarterly sheets).
# 2) Insert a new top data row into "总计" for the 2022-Q1 summary
#    (count=8, value=0.43), shifting the existing quarters down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New "2022-Q1" worksheet with the fund-holding detail table.
#    Duplicate an existing quarterly sheet so the header row/column-A index
#    style/sheetPr/pageMargins all come along for free, then swap in the
#    2022-Q1 numbers.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$template.Copy($totalSheet)

$ws = $wb.Worksheets.Item("2021-Q4 (2)")
$ws.Name = "2022-Q1"
$ws.Range("A1:H100").ClearContents()

# Header row (identical text/style to the other quarterly sheets)
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$data = @(
    @(0, "004871", "中银金融地产混合A",                          "3.67", "92.80", "3.39", "0.1244", 9),
    @(1, "001518", "万家瑞兴灵活配置混合",                        "3.16", "81.47", "3.21", "0.1014", 5),
    @(2, "011355", "华泰柏瑞港股通时代机遇混合型证券投资基金A",    "1.13", "90.93", "8.09", "0.0914", 2),
    @(3, "004099", "前海开源沪港深景气行业精选灵活配置混合",      "0.41", "93.07", "8.54", "0.0350", 5),
    @(4, "011356", "华泰柏瑞港股通时代机遇混合型证券投资基金C",    "0.40", "90.93", "8.09", "0.0324", 2),
    @(5, "010312", "中银金融地产混合C",                          "0.67", "92.80", "3.39", "0.0227", 9),
    @(6, "012315", "创金合信港股通成长股票型发起式证券投资基金A",  "0.19", "83.49", "7.93", "0.0151", 3),
    @(7, "012316", "创金合信港股通成长股票型发起式证券投资基金C",  "0.10", "83.49", "7.93", "0.0079", 3)
)

# Style source for the column-A index cells (rows 2..4 already carry it from
# the duplicated template; rows 5..9 are brand new and need it copied over).
$colAStyleSrc = $ws.Range("A2")
$colAStyleSrc.Copy()
$ws.Range("A5:A9").PasteSpecial(-4122)  # xlPasteFormats

$row = 2
foreach ($item in $data) {
    $ws.Range("A$row").Value = $item[0]

    # Text columns: force "@" text format before writing so values such as
    # "004871" or "3.67" are not coerced into numbers, then drop the format
    # again so the cell ends up with no explicit style (matching the other
    # detail sheets).
    $textRange = $ws.Range("B$row`:G$row")
    $textRange.NumberFormat = "@"
    $ws.Range("B$row").Value = $item[1]
    $ws.Range("C$row").Value = $item[2]
    $ws.Range("D$row").Value = $item[3]
    $ws.Range("E$row").Value = $item[4]
    $ws.Range("F$row").Value = $item[5]
    $ws.Range("G$row").Value = $item[6]
    $textRange.ClearFormats()

    $ws.Range("H$row").Value = $item[7]

    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 2) New top row in "总计" for the 2022-Q1 summary
# ---------------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")
$totalWs.Rows.Item(2).Insert()

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 8
$totalWs.Range("D2").Value = 0.43

# The freshly-inserted row inherited formatting from the header row above it;
# the target data rows carry no explicit style, so drop B2:D2's formatting...
$totalWs.Range("B2:D2").ClearFormats()
# ...and re-copy column-A styling from an existing data row (the insert also
# left A2 without the shared "index column" style).
$colASrc2 = $totalWs.Range("A3")
$colASrc2.Copy()
$totalWs.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

# Renumber the remaining index column (A3:A7 -> 1..5)
for ($r = 3; $r -le 7; $r++) {
    $totalWs.Range("A$r").Value = $r - 2
}

Write-Host "done"
